$wb = $excel.ActiveWorkbook

# Sheet 1
$ws = $wb.Worksheets.Item(1)
$ws.Range("H17").Value = 1287.0322
$ws.Range("I17").Value = 670.0
$ws.Range("J17").Value = 1378.4445
$ws.Range("K17").Value = 2010.0
$ws.Range("L17").Value = 4135.333500000001
$ws.Range("M17").Value = -1842.0
$ws.Range("N17").Value = -4471.333500000001
$ws.Range("H19").Value = 9433.143
$ws.Range("I19").Value = 9948.5
$ws.Range("K19").Value = 9948.5
$ws.Range("M19").Value = -9773.5
$ws.Range("H33").Value = 1588.1818
$ws.Range("I33").Value = 1719.0
$ws.Range("K33").Value = 1719.0
$ws.Range("M33").Value = -1490.0
$ws.Range("H62").Value = 16516.295
$ws.Range("I62").Value = 11393.546
$ws.Range("J62").Value = 25908.0
$ws.Range("K62").Value = 11393.546
$ws.Range("L62").Value = 25908.0
$ws.Range("M62").Value = -10769.546
$ws.Range("N62").Value = -27156.0
$ws.Range("H65").Value = 16516.295
$ws.Range("I65").Value = 11393.546
$ws.Range("J65").Value = 25908.0
$ws.Range("K65").Value = 56967.73
$ws.Range("L65").Value = 129540.0
$ws.Range("M65").Value = -53847.73
$ws.Range("N65").Value = -135780.0
$ws.Range("I74").Value = 250003120.0
$ws.Range("J74").Value = 14156.125
$ws.Range("K74").Value = 250003120.0
$ws.Range("L74").Value = 14156.125
$ws.Range("M74").Value = -250002184.0
$ws.Range("N74").Value = -16028.125
$ws.Range("I76").Value = 255416.0
$ws.Range("J76").Value = 62504500.0
$ws.Range("K76").Value = 255416.0
$ws.Range("L76").Value = 62504500.0
$ws.Range("M76").Value = -255101.0
$ws.Range("N76").Value = -62505130.0
$ws.Range("I77").Value = 250003120.0
$ws.Range("J77").Value = 14156.125
$ws.Range("K77").Value = 1250015600.0
$ws.Range("L77").Value = 70780.625
$ws.Range("M77").Value = -1250010920.0
$ws.Range("N77").Value = -80140.625
$ws.Range("I79").Value = 255416.0
$ws.Range("J79").Value = 62504500.0
$ws.Range("K79").Value = 255416.0
$ws.Range("L79").Value = 62504500.0
$ws.Range("M79").Value = -254324.0
$ws.Range("N79").Value = -62506684.0
$ws.Range("H86").Value = 51709304.0
$ws.Range("I86").Value = 65112964.0
$ws.Range("K86").Value = 65112964.0
$ws.Range("M86").Value = -65111841.0
$ws.Range("H89").Value = 51709304.0
$ws.Range("I89").Value = 65112964.0
$ws.Range("K89").Value = 325564820.0
$ws.Range("M89").Value = -325559204.0
$ws.Range("H92").Value = 169.88889
$ws.Range("I92").Value = 175.82353
$ws.Range("K92").Value = 175.82353
$ws.Range("M92").Value = 1072.17647
$ws.Range("H98").Value = 27780272.0
$ws.Range("I98").Value = 29414316.0
$ws.Range("K98").Value = 29414316.0
$ws.Range("M98").Value = -29412818.0
$ws.Range("H100").Value = 1092.7778
$ws.Range("I100").Value = 637.61536
$ws.Range("J100").Value = 2276.2
$ws.Range("K100").Value = 637.61536
$ws.Range("L100").Value = 2276.2
$ws.Range("M100").Value = -96.61536000000001
$ws.Range("N100").Value = -3358.2
$ws.Range("H103").Value = 927.6316
$ws.Range("I103").Value = 511.0
$ws.Range("J103").Value = 1230.6364
$ws.Range("K103").Value = 1533.0
$ws.Range("L103").Value = 3691.9092
$ws.Range("M103").Value = -947.0
$ws.Range("N103").Value = -4863.9092
$ws.Range("H106").Value = 50002236.0
$ws.Range("I106").Value = 62501810.0
$ws.Range("J106").Value = 3952.0
$ws.Range("K106").Value = 62501810.0
$ws.Range("L106").Value = 3952.0
$ws.Range("M106").Value = -62501179.0
$ws.Range("N106").Value = -5214.0
$ws.Range("H113").Value = 65228710.0
$ws.Range("I113").Value = 2043.75
$ws.Range("K113").Value = 2043.75
$ws.Range("M113").Value = 1210.25
$ws.Range("H116").Value = 13897979.0
$ws.Range("I116").Value = 22733640.0
$ws.Range("K116").Value = 22733640.0
$ws.Range("M116").Value = -22730198.0
$ws.Range("H122").Value = 27780272.0
$ws.Range("I122").Value = 29414316.0
$ws.Range("K122").Value = 88242948.0
$ws.Range("M122").Value = -88240498.0
$ws.Range("H131").Value = 991.25
$ws.Range("I131").Value = 991.25
$ws.Range("K131").Value = 2973.75
$ws.Range("M131").Value = 2066.25
$ws.Range("H132").Value = 0.0
$ws.Range("I132").Value = 0.0
$ws.Range("K132").Value = 0.0
$ws.Range("M132").Value = $null
$ws.Range("H137").Value = 8054.5
$ws.Range("I137").Value = 6973.6
$ws.Range("J137").Value = 8655.0
$ws.Range("K137").Value = 20920.8
$ws.Range("L137").Value = 25965.0
$ws.Range("M137").Value = -18370.8
$ws.Range("N137").Value = -31065.0
$ws.Range("H138").Value = 1729042.9
$ws.Range("I138").Value = 3098.0
$ws.Range("J138").Value = 2132987.5
$ws.Range("K138").Value = 9294.0
$ws.Range("L138").Value = 6398962.5
$ws.Range("M138").Value = -4154.0
$ws.Range("N138").Value = -6409242.5
$ws.Range("H141").Value = 2846.75
$ws.Range("I141").Value = 1000.0
$ws.Range("J141").Value = 3110.5715
$ws.Range("K141").Value = 3000.0
$ws.Range("L141").Value = 9331.7145
$ws.Range("M141").Value = 2180.0
$ws.Range("N141").Value = -19691.7145

# Sheet 2
$ws = $wb.Worksheets.Item(2)
$ws.Range("H2").Value = 71432056.0
$ws.Range("I2").Value = 4142.143
$ws.Range("J2").Value = 142859970.0
$ws.Range("K2").Value = 4142.143
$ws.Range("L2").Value = 142859970.0
$ws.Range("M2").Value = -4029.143
$ws.Range("N2").Value = -142860196.0
$ws.Range("H32").Value = 1304084.1
$ws.Range("I32").Value = 1304084.1
$ws.Range("K32").Value = 1304084.1
$ws.Range("M32").Value = -1303797.1
$ws.Range("H45").Value = 2874.0
$ws.Range("J45").Value = 2874.0
$ws.Range("L45").Value = 2874.0
$ws.Range("N45").Value = -3628.0
$ws.Range("H61").Value = 9885.134
$ws.Range("I61").Value = 2302.3333
$ws.Range("J61").Value = 11780.833
$ws.Range("K61").Value = 2302.3333
$ws.Range("L61").Value = 11780.833
$ws.Range("M61").Value = -2090.3333
$ws.Range("N61").Value = -12204.833
$ws.Range("H64").Value = 42999.145
$ws.Range("J64").Value = 42999.145
$ws.Range("L64").Value = 42999.145
$ws.Range("N64").Value = -43495.145
$ws.Range("H67").Value = 42999.145
$ws.Range("J67").Value = 42999.145
$ws.Range("L67").Value = 42999.145
$ws.Range("N67").Value = -44715.145
$ws.Range("H74").Value = 55765.535
$ws.Range("I74").Value = 77235.81
$ws.Range("J74").Value = 5668.222
$ws.Range("K74").Value = 77235.81
$ws.Range("L74").Value = 5668.222
$ws.Range("M74").Value = -76361.81
$ws.Range("N74").Value = -7416.222
$ws.Range("H77").Value = 55765.535
$ws.Range("I77").Value = 77235.81
$ws.Range("J77").Value = 5668.222
$ws.Range("K77").Value = 386179.05
$ws.Range("L77").Value = 28341.11
$ws.Range("M77").Value = -381811.05
$ws.Range("N77").Value = -37077.11
$ws.Range("H82").Value = 23396.75
$ws.Range("J82").Value = 23396.75
$ws.Range("L82").Value = 23396.75
$ws.Range("N82").Value = -24118.75
$ws.Range("H85").Value = 23396.75
$ws.Range("J85").Value = 23396.75
$ws.Range("L85").Value = 23396.75
$ws.Range("N85").Value = -25892.75
$ws.Range("H92").Value = 45873.0
$ws.Range("J92").Value = 45873.0
$ws.Range("L92").Value = 45873.0
$ws.Range("N92").Value = -50865.0
$ws.Range("H102").Value = 1585.6
$ws.Range("I102").Value = 1330.45
$ws.Range("J102").Value = 2095.9
$ws.Range("K102").Value = 1330.45
$ws.Range("L102").Value = 2095.9
$ws.Range("M102").Value = 291.55
$ws.Range("N102").Value = -5339.9
$ws.Range("H113").Value = 56272.25
$ws.Range("J113").Value = 56272.25
$ws.Range("L113").Value = 56272.25
$ws.Range("N113").Value = -64950.25
$ws.Range("H116").Value = 71432056.0
$ws.Range("I116").Value = 4142.143
$ws.Range("J116").Value = 142859970.0
$ws.Range("K116").Value = 4142.143
$ws.Range("L116").Value = 142859970.0
$ws.Range("M116").Value = -1848.143
$ws.Range("N116").Value = -142864558.0
$ws.Range("H122").Value = 25701.2
$ws.Range("I122").Value = 37502.0
$ws.Range("K122").Value = 112506.0
$ws.Range("M122").Value = -110056.0
$ws.Range("H132").Value = 4651.5894
$ws.Range("I132").Value = 1744.3334
$ws.Range("K132").Value = 5233.0002
$ws.Range("M132").Value = -2703.0002
$ws.Range("H136").Value = 9885.134
$ws.Range("I136").Value = 2302.3333
$ws.Range("J136").Value = 11780.833
$ws.Range("K136").Value = 6906.999899999999
$ws.Range("L136").Value = 35342.499
$ws.Range("M136").Value = -4356.999899999999
$ws.Range("N136").Value = -40442.499

# Sheet 3
$ws = $wb.Worksheets.Item(3)
$ws.Range("H3").Value = 71432056.0
$ws.Range("I3").Value = 4142.143
$ws.Range("J3").Value = 142859970.0
$ws.Range("K3").Value = 4142.143
$ws.Range("L3").Value = 142859970.0
$ws.Range("M3").Value = -4028.143
$ws.Range("N3").Value = -142860198.0
$ws.Range("H20").Value = 13892556.0
$ws.Range("I20").Value = 27780694.0
$ws.Range("K20").Value = 27780694.0
$ws.Range("M20").Value = -27780447.0
$ws.Range("H62").Value = 0.0
$ws.Range("J62").Value = 0.0
$ws.Range("L62").Value = 0.0
$ws.Range("N62").Value = $null
$ws.Range("H65").Value = 0.0
$ws.Range("J65").Value = 0.0
$ws.Range("L65").Value = 0.0
$ws.Range("N65").Value = $null
$ws.Range("H81").Value = 59750.11
$ws.Range("J81").Value = 59750.11
$ws.Range("L81").Value = 59750.11
$ws.Range("N81").Value = -61872.11
$ws.Range("H84").Value = 59750.11
$ws.Range("J84").Value = 59750.11
$ws.Range("L84").Value = 179250.33
$ws.Range("N84").Value = -189858.33
$ws.Range("H86").Value = 34210.742
$ws.Range("J86").Value = 2780.0
$ws.Range("L86").Value = 2780.0
$ws.Range("N86").Value = -5026.0
$ws.Range("H89").Value = 34210.742
$ws.Range("J89").Value = 2780.0
$ws.Range("L89").Value = 13900.0
$ws.Range("N89").Value = -25132.0
$ws.Range("H94").Value = 1982.7858
$ws.Range("I94").Value = 1539.0
$ws.Range("K94").Value = 1539.0
$ws.Range("M94").Value = -1088.0
$ws.Range("H99").Value = 22730522.0
$ws.Range("I99").Value = 5500.0
$ws.Range("K99").Value = 5500.0
$ws.Range("M99").Value = -4002.0
$ws.Range("H107").Value = 70379830.0
$ws.Range("I107").Value = 86619880.0
$ws.Range("K107").Value = 86619880.0
$ws.Range("M107").Value = -86617960.0
$ws.Range("H110").Value = 57687.5
$ws.Range("J110").Value = 57687.5
$ws.Range("L110").Value = 57687.5
$ws.Range("N110").Value = -65867.5
$ws.Range("H111").Value = 56458.668
$ws.Range("J111").Value = 56458.668
$ws.Range("L111").Value = 56458.668
$ws.Range("N111").Value = -64638.668
$ws.Range("H113").Value = 5495.125
$ws.Range("I113").Value = 5495.125
$ws.Range("K113").Value = 5495.125
$ws.Range("M113").Value = -3325.125
$ws.Range("H126").Value = 77283.5
$ws.Range("J126").Value = 77283.5
$ws.Range("L126").Value = 77283.5
$ws.Range("N126").Value = -87163.5
$ws.Range("H134").Value = 9623410.0
$ws.Range("I134").Value = 27780246.0
$ws.Range("J134").Value = 10967.588
$ws.Range("K134").Value = 83340738.0
$ws.Range("L134").Value = 32902.764
$ws.Range("M134").Value = -83338203.0
$ws.Range("N134").Value = -37972.764

# Sheet 4
$ws = $wb.Worksheets.Item(4)
$ws.Range("H31").Value = 8808.3545
$ws.Range("I31").Value = 2946.3076
$ws.Range("J31").Value = 13042.056
$ws.Range("K31").Value = 2946.3076
$ws.Range("L31").Value = 13042.056
$ws.Range("M31").Value = -2651.3076
$ws.Range("N31").Value = -13632.056
$ws.Range("H34").Value = 8808.3545
$ws.Range("I34").Value = 2946.3076
$ws.Range("J34").Value = 13042.056
$ws.Range("K34").Value = 2946.3076
$ws.Range("L34").Value = 13042.056
$ws.Range("M34").Value = -2744.3076
$ws.Range("N34").Value = -13446.056
$ws.Range("H58").Value = 8360.322
$ws.Range("I58").Value = 3928.8572
$ws.Range("J58").Value = 9652.833
$ws.Range("K58").Value = 3928.8572
$ws.Range("L58").Value = 9652.833
$ws.Range("M58").Value = -3725.8572
$ws.Range("N58").Value = -10058.833
$ws.Range("H62").Value = 3579.3
$ws.Range("I62").Value = 3356.8572
$ws.Range("K62").Value = 3356.8572
$ws.Range("M62").Value = -2732.8572
$ws.Range("H65").Value = 3579.3
$ws.Range("I65").Value = 3356.8572
$ws.Range("K65").Value = 16784.286
$ws.Range("M65").Value = -13664.286
$ws.Range("H99").Value = 6903.8076
$ws.Range("I99").Value = 7029.7646
$ws.Range("K99").Value = 7029.7646
$ws.Range("M99").Value = -5531.7646
$ws.Range("H111").Value = 90950.0
$ws.Range("J111").Value = 90950.0
$ws.Range("L111").Value = 90950.0
$ws.Range("N111").Value = -99130.0
$ws.Range("H112").Value = 92000.0
$ws.Range("J112").Value = 92000.0
$ws.Range("L112").Value = 92000.0
$ws.Range("N112").Value = -94954.0
$ws.Range("H122").Value = 4401.3335
$ws.Range("I122").Value = 3759.8333
$ws.Range("K122").Value = 11279.4999
$ws.Range("M122").Value = -8829.499899999999
$ws.Range("H126").Value = 6903.8076
$ws.Range("I126").Value = 7029.7646
$ws.Range("K126").Value = 21089.2938
$ws.Range("M126").Value = -18619.2938
$ws.Range("H132").Value = 6355.1816
$ws.Range("I132").Value = 4632.4546
$ws.Range("K132").Value = 13897.3638
$ws.Range("M132").Value = -11367.3638
$ws.Range("H134").Value = 7265.0
$ws.Range("I134").Value = 2620.375
$ws.Range("J134").Value = 9329.277
$ws.Range("K134").Value = 7861.125
$ws.Range("L134").Value = 27987.831
$ws.Range("M134").Value = -5326.125
$ws.Range("N134").Value = -33057.831
$ws.Range("H136").Value = 8360.322
$ws.Range("I136").Value = 3928.8572
$ws.Range("J136").Value = 9652.833
$ws.Range("K136").Value = 11786.5716
$ws.Range("L136").Value = 28958.499
$ws.Range("M136").Value = -9236.5716
$ws.Range("N136").Value = -34058.499

# Sheet 5
$ws = $wb.Worksheets.Item(5)
$ws.Range("H2").Value = 71764.4
$ws.Range("I2").Value = 11164.852
$ws.Range("J2").Value = 276287.88
$ws.Range("K2").Value = 66989.11200000001
$ws.Range("L2").Value = 1657727.28
$ws.Range("M2").Value = -66876.11200000001
$ws.Range("N2").Value = -1657953.28
$ws.Range("H5").Value = 3046.7368
$ws.Range("I5").Value = 671.8182
$ws.Range("J5").Value = 6312.25
$ws.Range("K5").Value = 2015.4546
$ws.Range("L5").Value = 18936.75
$ws.Range("M5").Value = -1903.4546
$ws.Range("N5").Value = -19160.75
$ws.Range("H26").Value = 282.89474
$ws.Range("I26").Value = 112.0
$ws.Range("J26").Value = 328.46667
$ws.Range("K26").Value = 336.0
$ws.Range("L26").Value = 985.4000100000001
$ws.Range("M26").Value = -48.0
$ws.Range("N26").Value = -1561.40001
$ws.Range("H32").Value = 100.0
$ws.Range("I32").Value = 0.0
$ws.Range("K32").Value = 0.0
$ws.Range("M32").Value = $null
$ws.Range("I38").Value = 24.285715
$ws.Range("K38").Value = 72.857145
$ws.Range("M38").Value = 274.142855
$ws.Range("H40").Value = 292.5
$ws.Range("I40").Value = 220.0
$ws.Range("J40").Value = 316.66666
$ws.Range("K40").Value = 880.0
$ws.Range("L40").Value = 1266.66664
$ws.Range("M40").Value = -811.0
$ws.Range("N40").Value = -1404.66664
$ws.Range("H55").Value = 6258807.0
$ws.Range("I55").Value = 0.0
$ws.Range("K55").Value = 0.0
$ws.Range("M55").Value = $null
$ws.Range("H68").Value = 22225480.0
$ws.Range("I68").Value = 50000652.0
$ws.Range("J68").Value = 5344.6
$ws.Range("K68").Value = 150001956.0
$ws.Range("L68").Value = 16033.8
$ws.Range("M68").Value = -150001145.0
$ws.Range("N68").Value = -17655.8
$ws.Range("H71").Value = 22225480.0
$ws.Range("I71").Value = 50000652.0
$ws.Range("J71").Value = 5344.6
$ws.Range("K71").Value = 450005868.0
$ws.Range("L71").Value = 48101.4
$ws.Range("M71").Value = -450001812.0
$ws.Range("N71").Value = -56213.4
$ws.Range("H107").Value = 7692971.0
$ws.Range("I107").Value = 265.70587
$ws.Range("J107").Value = 22223636.0
$ws.Range("K107").Value = 797.11761
$ws.Range("L107").Value = 66670908.0
$ws.Range("M107").Value = 1122.88239
$ws.Range("N107").Value = -66674748.0
$ws.Range("H113").Value = 2655.2222
$ws.Range("J113").Value = 3271.1428
$ws.Range("L113").Value = 9813.4284
$ws.Range("N113").Value = -14153.4284
$ws.Range("H129").Value = 18575192.0
$ws.Range("I129").Value = 939.625
$ws.Range("J129").Value = 33434594.0
$ws.Range("K129").Value = 2818.875
$ws.Range("L129").Value = 100303782.0
$ws.Range("M129").Value = 2181.125
$ws.Range("N129").Value = -100313782.0
$ws.Range("H131").Value = 2250.3403
$ws.Range("I131").Value = 1304.0
$ws.Range("K131").Value = 3912.0
$ws.Range("M131").Value = 1128.0
$ws.Range("H132").Value = 6241.517
$ws.Range("J132").Value = 12049.833
$ws.Range("L132").Value = 108448.497
$ws.Range("N132").Value = -113508.497
$ws.Range("H135").Value = 3046.7368
$ws.Range("I135").Value = 671.8182
$ws.Range("J135").Value = 6312.25
$ws.Range("K135").Value = 6046.3638
$ws.Range("L135").Value = 56810.25
$ws.Range("M135").Value = -3511.3638
$ws.Range("N135").Value = -61880.25
$ws.Range("H139").Value = 54267.8
$ws.Range("I139").Value = 62077.824
$ws.Range("K139").Value = 186233.472
$ws.Range("M139").Value = -181093.472

# Sheet 6
$ws = $wb.Worksheets.Item(6)
$ws.Range("H80").Value = 203089.4
$ws.Range("I80").Value = 2815.6667
$ws.Range("J80").Value = 503500.0
$ws.Range("K80").Value = 2815.6667
$ws.Range("L80").Value = 503500.0
$ws.Range("M80").Value = -1817.6667
$ws.Range("N80").Value = -505496.0
$ws.Range("H83").Value = 203089.4
$ws.Range("I83").Value = 2815.6667
$ws.Range("J83").Value = 503500.0
$ws.Range("K83").Value = 14078.3335
$ws.Range("L83").Value = 2517500.0
$ws.Range("M83").Value = -9086.3335
$ws.Range("N83").Value = -2527484.0
$ws.Range("H95").Value = 53921.0
$ws.Range("J95").Value = 53921.0
$ws.Range("L95").Value = 53921.0
$ws.Range("N95").Value = -59413.0
$ws.Range("H96").Value = 50935.25
$ws.Range("J96").Value = 50935.25
$ws.Range("L96").Value = 50935.25
$ws.Range("N96").Value = -56427.25
$ws.Range("H97").Value = 3053.3333
$ws.Range("I97").Value = 2924.5
$ws.Range("J97").Value = 3311.0
$ws.Range("K97").Value = 2924.5
$ws.Range("L97").Value = 3311.0
$ws.Range("M97").Value = -2428.5
$ws.Range("N97").Value = -4303.0
$ws.Range("H113").Value = 5653.447
$ws.Range("I113").Value = 3095.8333
$ws.Range("K113").Value = 3095.8333
$ws.Range("M113").Value = -925.8332999999998
$ws.Range("H122").Value = 1908483.8
$ws.Range("I122").Value = 2684737.8
$ws.Range("J122").Value = 3132.9092
$ws.Range("K122").Value = 8054213.399999999
$ws.Range("L122").Value = 9398.7276
$ws.Range("M122").Value = -8051763.399999999
$ws.Range("N122").Value = -14298.7276
$ws.Range("H126").Value = 6290.706
$ws.Range("I126").Value = 2459.5715
$ws.Range("J126").Value = 8972.5
$ws.Range("K126").Value = 7378.7145
$ws.Range("L126").Value = 26917.5
$ws.Range("M126").Value = -4908.7145
$ws.Range("N126").Value = -31857.5
$ws.Range("H132").Value = 2308.3333
$ws.Range("I132").Value = 2304.7036
$ws.Range("J132").Value = 2324.6667
$ws.Range("K132").Value = 6914.110799999999
$ws.Range("L132").Value = 6974.000100000001
$ws.Range("M132").Value = -4384.110799999999
$ws.Range("N132").Value = -12034.0001

# Sheet 7
$ws = $wb.Worksheets.Item(7)
$ws.Range("H7").Value = 5349.5654
$ws.Range("I7").Value = 2942.75
$ws.Range("K7").Value = 2942.75
$ws.Range("M7").Value = -2830.75
$ws.Range("H16").Value = 1438.2
$ws.Range("I16").Value = 1472.5
$ws.Range("J16").Value = 1301.0
$ws.Range("K16").Value = 1472.5
$ws.Range("L16").Value = 1301.0
$ws.Range("M16").Value = -1302.5
$ws.Range("N16").Value = -1641.0
$ws.Range("H22").Value = 6212.9
$ws.Range("I22").Value = 645.6667
$ws.Range("K22").Value = 645.6667
$ws.Range("M22").Value = -350.6667
$ws.Range("H27").Value = 6212.9
$ws.Range("I27").Value = 645.6667
$ws.Range("K27").Value = 645.6667
$ws.Range("M27").Value = -538.6667
$ws.Range("H46").Value = 4632280.5
$ws.Range("I46").Value = 1876.4445
$ws.Range("J46").Value = 7410523.5
$ws.Range("K46").Value = 1876.4445
$ws.Range("L46").Value = 7410523.5
$ws.Range("M46").Value = -1688.4445
$ws.Range("N46").Value = -7410899.5
$ws.Range("H61").Value = 4580.0
$ws.Range("I61").Value = 2772.25
$ws.Range("K61").Value = 2772.25
$ws.Range("M61").Value = -2570.25
$ws.Range("H82").Value = 387377.3
$ws.Range("I82").Value = 669027.6
$ws.Range("K82").Value = 669027.6
$ws.Range("M82").Value = -668666.6
$ws.Range("H85").Value = 387377.3
$ws.Range("I85").Value = 669027.6
$ws.Range("K85").Value = 669027.6
$ws.Range("M85").Value = -667779.6
$ws.Range("H93").Value = 2836.0908
$ws.Range("I93").Value = 2724.875
$ws.Range("J93").Value = 3132.6667
$ws.Range("K93").Value = 2724.875
$ws.Range("L93").Value = 3132.6667
$ws.Range("M93").Value = -1476.875
$ws.Range("N93").Value = -5628.6667
$ws.Range("H100").Value = 4409.4546
$ws.Range("I100").Value = 2043.0
$ws.Range("K100").Value = 2043.0
$ws.Range("M100").Value = -1502.0
$ws.Range("H113").Value = 4580.0
$ws.Range("I113").Value = 2772.25
$ws.Range("K113").Value = 2772.25
$ws.Range("M113").Value = -602.25
$ws.Range("H114").Value = 55938.0
$ws.Range("J114").Value = 55938.0
$ws.Range("L114").Value = 55938.0
$ws.Range("N114").Value = -64616.0
$ws.Range("H122").Value = 3989.8262
$ws.Range("I122").Value = 2692.1177
$ws.Range("J122").Value = 7666.6665
$ws.Range("K122").Value = 8076.353099999999
$ws.Range("L122").Value = 22999.9995
$ws.Range("M122").Value = -5626.353099999999
$ws.Range("N122").Value = -27899.9995
$ws.Range("H125").Value = 0.0
$ws.Range("J125").Value = 0.0
$ws.Range("L125").Value = 0.0
$ws.Range("N125").Value = $null
$ws.Range("H126").Value = 5349.5654
$ws.Range("I126").Value = 2942.75
$ws.Range("K126").Value = 8828.25
$ws.Range("M126").Value = -6358.25
$ws.Range("H136").Value = 7176.426
$ws.Range("I136").Value = 2861.04
$ws.Range("J136").Value = 10896.586
$ws.Range("K136").Value = 8583.119999999999
$ws.Range("L136").Value = 32689.758
$ws.Range("M136").Value = -6033.119999999999
$ws.Range("N136").Value = -37789.758

# Sheet 8
$ws = $wb.Worksheets.Item(8)
$ws.Range("H8").Value = 0.0
$ws.Range("I8").Value = 0.0
$ws.Range("K8").Value = 0.0
$ws.Range("M8").Value = $null
$ws.Range("H14").Value = 200200600.0
$ws.Range("I14").Value = 200200600.0
$ws.Range("K14").Value = 200200600.0
$ws.Range("M14").Value = -200200432.0
$ws.Range("H69").Value = 11000.0
$ws.Range("I69").Value = 11000.0
$ws.Range("K69").Value = 11000.0
$ws.Range("M69").Value = -10251.0
$ws.Range("H72").Value = 11000.0
$ws.Range("I72").Value = 11000.0
$ws.Range("K72").Value = 33000.0
$ws.Range("M72").Value = -29256.0
$ws.Range("H81").Value = 21057408.0
$ws.Range("I81").Value = 5199999.0
$ws.Range("J81").Value = 25021760.0
$ws.Range("K81").Value = 10399998.0
$ws.Range("L81").Value = 50043520.0
$ws.Range("M81").Value = -10398937.0
$ws.Range("N81").Value = -50045642.0
$ws.Range("H84").Value = 21057408.0
$ws.Range("I84").Value = 5199999.0
$ws.Range("J84").Value = 25021760.0
$ws.Range("K84").Value = 51999990.0
$ws.Range("L84").Value = 250217600.0
$ws.Range("M84").Value = -51994686.0
$ws.Range("N84").Value = -250228208.0
$ws.Range("H114").Value = 56546.0
$ws.Range("J114").Value = 56546.0
$ws.Range("L114").Value = 56546.0
$ws.Range("N114").Value = -65224.0
$ws.Range("H122").Value = 10725980.0
$ws.Range("I122").Value = 21002162.0
$ws.Range("K122").Value = 63006486.0
$ws.Range("M122").Value = -63004036.0
$ws.Range("H126").Value = 1602.2424
$ws.Range("I126").Value = 1651.1305
$ws.Range("K126").Value = 4953.3915
$ws.Range("M126").Value = -2483.3915
$ws.Range("H132").Value = 25646246.0
$ws.Range("I132").Value = 30308002.0
$ws.Range("J132").Value = 6583.0
$ws.Range("K132").Value = 90924006.0
$ws.Range("L132").Value = 19749.0
$ws.Range("M132").Value = -90921476.0
$ws.Range("N132").Value = -24809.0
$ws.Range("H136").Value = 33669740.0
$ws.Range("I136").Value = 45455692.0
$ws.Range("J136").Value = 1258376.2
$ws.Range("K136").Value = 136367076.0
$ws.Range("L136").Value = 3775128.6
$ws.Range("M136").Value = -136364526.0
$ws.Range("N136").Value = -3780228.6
